$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 98
$ws.Range("A98").Value = 111790914
$ws.Range("B98").Value = 90662
$ws.Range("D98").Value = "LC"
$ws.Range("E98").Value = 4363
$ws.Range("F98").Value = "Zontaggsvamp"
$ws.Range("G98").Value = "Hydnellum concrescens"
$ws.Range("H98").Value = "(Pers.) Banker"
$ws.Range("P98").Value = "Storvreta (Storvreta), Upl"
$ws.Range("Q98").Value = 650131
$ws.Range("R98").Value = 6648444
$ws.Range("S98").Value = 100
$ws.Range("Z98").Value = "'18:43"
$ws.Range("AB98").Value = "'18:43"
$ws.Range("AC98").Value = $null

# Row 99
$ws.Range("A99").Value = 111789319
$ws.Range("B99").Value = 88915
$ws.Range("D99").Value = "NT"
$ws.Range("E99").Value = 5734
$ws.Range("F99").Value = "Druvfingersvamp"
$ws.Range("G99").Value = "Ramaria botrytis"
$ws.Range("H99").Value = "(Pers.:Fr.) Bourdot"
$ws.Range("Q99").Value = 650042
$ws.Range("R99").Value = 6648755
$ws.Range("Z99").Value = "'16:58"
$ws.Range("AB99").Value = "'16:58"
$ws.Range("AC99").Value = "Mitt i en stigkorsning"

# Row 100
$ws.Range("Q100").Value = 649875
$ws.Range("R100").Value = 6648703

# Row 101
$ws.Range("A101").Value = 111790472
$ws.Range("P101").Value = "Storvreta, Upl"
$ws.Range("Q101").Value = 650061
$ws.Range("R101").Value = 6648732
$ws.Range("S101").Value = 25
$ws.Range("Z101").Value = "'18:09"
$ws.Range("AB101").Value = "'18:09"
$ws.Range("AC101").Value = "Precis på/vid stig. Andra mycel åt båda håll."

# Row 102
$ws.Range("A102").Value = 111791986
$ws.Range("B102").Value = 88982
$ws.Range("D102").Value = "NT"
$ws.Range("E102").Value = 937
$ws.Range("F102").Value = "Vit vedfingersvamp"
$ws.Range("G102").Value = "Lentaria epichnoa"
$ws.Range("H102").Value = "(Fr.) Corner"
$ws.Range("P102").Value = "Storvreta, Upl"
$ws.Range("Q102").Value = 650061
$ws.Range("R102").Value = 6648732
$ws.Range("S102").Value = 25
$ws.Range("Y102").Value = "'2022-10-01"
$ws.Range("Z102").Value = "'19:25"
$ws.Range("AA102").Value = "'2022-10-01"
$ws.Range("AB102").Value = "'19:25"
$ws.Range("AC102").Value = "Noterad ett år sent efter att ha problem med att rapportera! Två kraftigt rötade asplågor intill stig löpande N->S mellan Himmelsvägen och hästgården/huset efter skogen. Delar lågor med biskopsmössor."

# Row 103
$ws.Range("A103").Value = 111790550
$ws.Range("B103").Value = 90655
$ws.Range("D103").Value = "VU"
$ws.Range("E103").Value = 150
$ws.Range("F103").Value = "Grangråticka"
$ws.Range("G103").Value = "Boletopsis leucomelaena"
$ws.Range("H103").Value = "(Pers.) Fayod"
$ws.Range("Q103").Value = 649935
$ws.Range("R103").Value = 6648620
$ws.Range("S103").Value = 20
$ws.Range("Z103").Value = "'18:21"
$ws.Range("AB103").Value = "'18:21"

# Row 104
$ws.Range("A104").Value = 111789477
$ws.Range("B104").Value = 88909
$ws.Range("D104").Value = "VU"
$ws.Range("E104").Value = 720
$ws.Range("F104").Value = "Violgubbe"
$ws.Range("G104").Value = "Gomphus clavatus"
$ws.Range("H104").Value = "(Pers.) Gray"
$ws.Range("Q104").Value = 649983
$ws.Range("R104").Value = 6648731
$ws.Range("Z104").Value = "'17:09"
$ws.Range("AB104").Value = "'17:09"
$ws.Range("AC104").Value = "Mellan stenblock ca 50 m bort ifrån föregående rapporterade mycel. Bland ett myller av granlågor."
$ws.Range("AH104").Value = $null

# Row 105
$ws.Range("A105").Value = 111789368
$ws.Range("B105").Value = 88909
$ws.Range("D105").Value = "VU"
$ws.Range("E105").Value = 720
$ws.Range("F105").Value = "Violgubbe"
$ws.Range("G105").Value = "Gomphus clavatus"
$ws.Range("H105").Value = "(Pers.) Gray"
$ws.Range("P105").Value = "Storvreta (Storvreta), Upl"
$ws.Range("Q105").Value = 650001
$ws.Range("R105").Value = 6648759
$ws.Range("S105").Value = 10
$ws.Range("Z105").Value = "'17:09"
$ws.Range("AB105").Value = "'17:09"
$ws.Range("AC105").Value = "I ett ca 5 meter långt stråk med många fruktkroppar."
$ws.Range("AH105").Value = "Blåbärsgranskog"

# Row 106
$ws.Range("A106").Value = 111789261
$ws.Range("B106").Value = 88915
$ws.Range("E106").Value = 5734
$ws.Range("F106").Value = "Druvfingersvamp"
$ws.Range("G106").Value = "Ramaria botrytis"
$ws.Range("H106").Value = "(Pers.:Fr.) Bourdot"
$ws.Range("Q106").Value = 650061
$ws.Range("R106").Value = 6648732
$ws.Range("Y106").Value = "'2023-08-30"
$ws.Range("Z106").Value = "'17:02"
$ws.Range("AA106").Value = "'2023-08-30"
$ws.Range("AB106").Value = "'17:02"
$ws.Range("AC106").Value = $null

# Row 107
$ws.Range("Q107").Value = 650071
$ws.Range("R107").Value = 6648472

# Row 108
$ws.Range("Q108").Value = 650338
$ws.Range("R108").Value = 6648422
$ws.Range("Z108").Value = $null
$ws.Range("AB108").Value = $null

# Row 109
$ws.Range("A109").Value = 111984708
$ws.Range("B109").Value = 88909
$ws.Range("D109").Value = "VU"
$ws.Range("E109").Value = 720
$ws.Range("F109").Value = "Violgubbe"
$ws.Range("G109").Value = "Gomphus clavatus"
$ws.Range("H109").Value = "(Pers.) Gray"
$ws.Range("Q109").Value = 650057
$ws.Range("R109").Value = 6648630
$ws.Range("Z109").Value = "'14:54"
$ws.Range("AB109").Value = "'14:54"
$ws.Range("AC109").Value = "Under granar, precis intill kärret och en halvmeter från stig. Två samlingar."

# Row 110
$ws.Range("A110").Value = 111984394
$ws.Range("B110").Value = 90662
$ws.Range("D110").Value = "LC"
$ws.Range("E110").Value = 4363
$ws.Range("F110").Value = "Zontaggsvamp"
$ws.Range("G110").Value = "Hydnellum concrescens"
$ws.Range("H110").Value = "(Pers.) Banker"
$ws.Range("Q110").Value = 650057
$ws.Range("R110").Value = 6648630
$ws.Range("Z110").Value = "'13:16"
$ws.Range("AB110").Value = "'13:16"
$ws.Range("AC110").Value = "Nedanför granbacke i en stig (kärr ca 40 meter NV)"

# Row 111
$ws.Range("Q111").Value = 650049
$ws.Range("R111").Value = 6648645
$ws.Range("Z111").Value = $null
$ws.Range("AB111").Value = $null

# Row 112
$ws.Range("Q112").Value = 650012
$ws.Range("R112").Value = 6648763
$ws.Range("Z112").Value = $null
$ws.Range("AB112").Value = $null

# Row 113
$ws.Range("Q113").Value = 649886
$ws.Range("R113").Value = 6648971
$ws.Range("Z113").Value = $null
$ws.Range("AB113").Value = $null
